$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: G1/H1/I1 were "Didymin1"/"Didymin2"/"Didymin3" -> now "MA1"/"MA2"/"MA3"
$ws.Range("G1").Value = "MA1"
$ws.Range("H1").Value = "MA2"
$ws.Range("I1").Value = "MA3"

# Data edits on rows 64, 65, 68 (Model3 / HMDB0009211 / HMDB0134992 rows)
$ws.Range("B64").Value = 1
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0

$ws.Range("B65").Value = 1
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0

$ws.Range("B68").Value = 1
$ws.Range("D68").Value = 0

# View state: scroll to row 43, select B66
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("B66").Select()
